$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Review Records")
$ws.Activate()

# Row 18 (No = 17) - fill in ORM record details
$ws.Range("B18").Value = 14011917
$ws.Range("C18").Value = "TUCMS.docx"
$ws.Range("D18").Value = "Kaung Myat Bo"
$ws.Range("E17").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("E18").Value = 43479

# Row 19 (No = 18) - fill in ORM record details
$ws.Range("B19").Value = 14011918
$ws.Range("C19").Value = "HLD.docx"
$ws.Range("D19").Value = "Kaung Myat Bo"
$ws.Range("E17").Copy()
$ws.Range("E19").PasteSpecial(-4122)
$ws.Range("E19").Value = 43479

# Update the selected cell to match the saved view state
$ws.Range("C19").Select()

$wb.Save()
